# Apply the "750 last non-tied games" update.
#
# The workbook recomputes the 11 most common winning/losing scores (PF) from
# the last N non-tied regular-season games. N changes from 500 to 750, which
# changes the underlying "W ex R" / "L ex R" score-pair strings in L2:L12 and
# S2:S12. Every other changed cell in the sheet (G1, J1, Q28, X28, the M/N/O/P/Q
# and T/U/V/W/X helper columns, the data-table G2:G5/J2:J5, and the B9:C12
# label strings) is a formula that is derived from those inputs, so updating
# the 22 source strings plus the comment text and letting Excel recalculate
# reproduces the full change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the explanatory comment on L1 (500 -> 750, twice) ---------------
$comment = $ws.Range("L1").Comment
$newCommentText = "Die elf häufigsten W scores der letzten 750 non tied RS games`n`nresults |> filter(Result == ""W"", Week < 30) |> tail(750) |> pull(PF) |> table() |> sort(decreasing = T) |> t() |> t() |> head(11)"
$comment.Text($newCommentText)

# --- Update the "W ex R" source score pairs (L2:L12) -------------------------
$ws.Range("L2").Value = "  27   68"
$ws.Range("L3").Value = "  20   56"
$ws.Range("L4").Value = "  31   48"
$ws.Range("L5").Value = "  24   46"
$ws.Range("L6").Value = "  23   39"
$ws.Range("L7").Value = "  30   39"
$ws.Range("L8").Value = "  34   36"
$ws.Range("L9").Value = "  17   30"
$ws.Range("L10").Value = "  19   29"
$ws.Range("L11").Value = "  28   27"
$ws.Range("L12").Value = "  26   26"

# --- Update the "L ex R" source score pairs (S2:S12) -------------------------
$ws.Range("S2").Value = "  17   80"
$ws.Range("S3").Value = "  10   65"
$ws.Range("S4").Value = "  16   56"
$ws.Range("S5").Value = "  20   46"
$ws.Range("S6").Value = "  13   44"
$ws.Range("S7").Value = "  24   41"
$ws.Range("S8").Value = "  14   38"
$ws.Range("S9").Value = "  3    33"
$ws.Range("S10").Value = "  7    28"
$ws.Range("S11").Value = "  21   27"
$ws.Range("S12").Value = "  9    25"

# Force a full recalculation so every dependent formula (including the
# G2:G5 / J2:J5 What-If data tables) picks up the new source data.
$excel.CalculateFull()
